$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("worksheet_no2")

# Insert a new column before column A, shifting the existing data
# (5/13, 6/13, absolute_diff, relative_diff columns) one column to the right.
$ws.Range("A1:A6").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("A1").Value = "metric"

# New metric-name labels for each data row.
$ws.Range("A2").Value = "sessions"
$ws.Range("A3").Value = "transactions"
$ws.Range("A4").Value = "QTY"
$ws.Range("A5").Value = "ECR"
$ws.Range("A6").Value = "addsToCart"

# Rename the "month_number" shared header (used on worksheet_no1) to
# "year_after_2000".
$ws1 = $wb.Worksheets.Item("worksheet_no1")
$ws1.Range("B1").Value = "year_after_2000"
